$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New "Creator" column (E) with header + two "N.N." entries
$ws.Range("E1").Value = "Creator"
$ws.Range("E2").Value = "N.N."
$ws.Range("E4").Value = "N.N."

# Set the width of column E (closest achievable value to the authored 25.47)
$ws.Columns.Item(5).ColumnWidth = 24.65

# Move selection to E12 to match the author's final cursor position
$ws.Range("E12").Select()
